$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 3: Jane Gichohi
Set-TextValue $ws.Range("A3") "Jane Gichohi"
Set-TextValue $ws.Range("B3") "2.00"
Set-TextValue $ws.Range("D3") "-13.00"
Set-TextValue $ws.Range("E3") "13.33%"

# Row 4: Mirriam Makau
Set-TextValue $ws.Range("A4") "Mirriam Makau"
Set-TextValue $ws.Range("B4") "1.00"
Set-TextValue $ws.Range("C4") "12.00"
Set-TextValue $ws.Range("D4") "-11.00"
Set-TextValue $ws.Range("E4") "8.33%"

# Row 5: KD Totals
Set-TextValue $ws.Range("A5") "KD Totals"
Set-TextValue $ws.Range("B5") "5.00"
Set-TextValue $ws.Range("C5") "39.00"
Set-TextValue $ws.Range("D5") "-34.00"
Set-TextValue $ws.Range("E5") "12.78%"

# Remove old rows 6 and 7 (Victor Njogu, old KD Totals) - shift cells up
$ws.Range("A6:E7").Delete()
